$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: new deposito transaction
$ws.Range("B19").Value = "2025-03-08 12:32:14"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "74887540"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "deposito"
$ws.Range("E19").Value = 120

# Row 20: new retiro transaction
$ws.Range("B20").Value = "2025-03-08 12:34:31"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "74887540"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "retiro"
$ws.Range("E20").Value = 100
